$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
$ws1 = $wb.Worksheets.Item("Overview")
$ws1.Range("A2").Value = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws1.Range("A3").Value = "ffffa596e9e9-65cf-4444-a47c-43744ad452c8.md"

# --- zh-cn sheet ---
$ws2 = $wb.Worksheets.Item("zh-cn")
$ws2.Range("A2").Value = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws2.Range("C2").Value = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.zh-cn.xlf"
$ws2.Range("D2").Value = "2016-03-10 23:19:53"
$ws2.Range("E2").Value = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws2.Range("F2").Value = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.zh-cn.xlf"
$ws2.Range("G2").Value = "2016-03-10 23:20:25"

$ws2.Range("A3").Value = "ffffa596e9e9-65cf-4444-a47c-43744ad452c8.md"
$ws2.Range("C3").Value = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.zh-cn.xlf"
$ws2.Range("D3").Value = "2016-03-10 23:19:53"
$ws2.Range("E3").Value = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws2.Range("F3").Value = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.zh-cn.xlf"
$ws2.Range("G3").Value = "2016-03-10 23:20:25"

$ws2.Range("D4").Value = "0001-01-01 00:00:00"
$ws2.Range("G4").Value = "0001-01-01 00:00:00"

# --- de-de sheet ---
$ws3 = $wb.Worksheets.Item("de-de")
$ws3.Range("A2").Value = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws3.Range("C2").Value = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.de-de.xlf"
$ws3.Range("D2").Value = "2016-03-10 23:19:59"
$ws3.Range("E2").Value = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws3.Range("F2").Value = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.de-de.xlf"
$ws3.Range("G2").Value = "2016-03-10 23:20:42"

$ws3.Range("A3").Value = "ffffa596e9e9-65cf-4444-a47c-43744ad452c8.md"
$ws3.Range("C3").Value = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.de-de.xlf"
$ws3.Range("D3").Value = "2016-03-10 23:19:59"
$ws3.Range("E3").Value = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws3.Range("F3").Value = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.de-de.xlf"
$ws3.Range("G3").Value = "2016-03-10 23:20:42"

$ws3.Range("D4").Value = "0001-01-01 00:00:00"
$ws3.Range("G4").Value = "0001-01-01 00:00:00"

# --- Hyperlinks: update display text and target URLs to match new file names ---

# Overview sheet hyperlinks
$ws1.Hyperlinks.Item(1).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws1.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTest/oltest/blob/16e8fb4e99ce355c090d21ec564c8937f9d6d0d5/e2e/a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws1.Hyperlinks.Item(2).TextToDisplay = "ffffa596e9e9-65cf-4444-a47c-43744ad452c8.md"
$ws1.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTest/oltest/blob/16e8fb4e99ce355c090d21ec564c8937f9d6d0d5/e2e/ffffa596e9e9-65cf-4444-a47c-43744ad452c8.md"

# zh-cn sheet hyperlinks (A2,C2,E2,F2,A3,C3,E3,F3)
$ws2.Hyperlinks.Item(1).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws2.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTest/oltest/blob/16e8fb4e99ce355c090d21ec564c8937f9d6d0d5/e2e/a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws2.Hyperlinks.Item(2).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.zh-cn.xlf"
$ws2.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e28828a9c69a1ed50a56d3fa2be1d3bda7d5d48d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.zh-cn.xlf"
$ws2.Hyperlinks.Item(3).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws2.Hyperlinks.Item(3).Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/43346ac97fef70fa7a078145b50ec07efbb26d03/e2e/a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws2.Hyperlinks.Item(4).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.zh-cn.xlf"
$ws2.Hyperlinks.Item(4).Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/36ad4be5af9f914141ed1f606a0f9e0cbfb1cb9f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.zh-cn.xlf"
$ws2.Hyperlinks.Item(5).TextToDisplay = "ffffa596e9e9-65cf-4444-a47c-43744ad452c8.md"
$ws2.Hyperlinks.Item(5).Address = "https://github.com/OpenLocalizationTest/oltest/blob/16e8fb4e99ce355c090d21ec564c8937f9d6d0d5/e2e/ffffa596e9e9-65cf-4444-a47c-43744ad452c8.md"
$ws2.Hyperlinks.Item(6).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.zh-cn.xlf"
$ws2.Hyperlinks.Item(6).Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e28828a9c69a1ed50a56d3fa2be1d3bda7d5d48d/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.zh-cn.xlf"
$ws2.Hyperlinks.Item(7).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws2.Hyperlinks.Item(7).Address = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/43346ac97fef70fa7a078145b50ec07efbb26d03/e2e/a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws2.Hyperlinks.Item(8).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.zh-cn.xlf"
$ws2.Hyperlinks.Item(8).Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/36ad4be5af9f914141ed1f606a0f9e0cbfb1cb9f/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.zh-cn.xlf"

# de-de sheet hyperlinks (A2,C2,E2,F2,A3,C3,E3,F3)
$ws3.Hyperlinks.Item(1).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws3.Hyperlinks.Item(1).Address = "https://github.com/OpenLocalizationTest/oltest/blob/16e8fb4e99ce355c090d21ec564c8937f9d6d0d5/e2e/a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws3.Hyperlinks.Item(2).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.de-de.xlf"
$ws3.Hyperlinks.Item(2).Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/483e07b192e98807fbf5d1bd3b5792b4fc7706fb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.de-de.xlf"
$ws3.Hyperlinks.Item(3).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws3.Hyperlinks.Item(3).Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/50d8f0ef45a01b6758b64a615b34fe57373bc2b2/e2e/a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws3.Hyperlinks.Item(4).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.de-de.xlf"
$ws3.Hyperlinks.Item(4).Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e6bf38c351c2af39466f1fd9fe06bae116651949/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.de-de.xlf"
$ws3.Hyperlinks.Item(5).TextToDisplay = "ffffa596e9e9-65cf-4444-a47c-43744ad452c8.md"
$ws3.Hyperlinks.Item(5).Address = "https://github.com/OpenLocalizationTest/oltest/blob/16e8fb4e99ce355c090d21ec564c8937f9d6d0d5/e2e/ffffa596e9e9-65cf-4444-a47c-43744ad452c8.md"
$ws3.Hyperlinks.Item(6).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.de-de.xlf"
$ws3.Hyperlinks.Item(6).Address = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/483e07b192e98807fbf5d1bd3b5792b4fc7706fb/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.de-de.xlf"
$ws3.Hyperlinks.Item(7).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws3.Hyperlinks.Item(7).Address = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/50d8f0ef45a01b6758b64a615b34fe57373bc2b2/e2e/a8b88462-4116-4a82-a657-005344fb7aed.md"
$ws3.Hyperlinks.Item(8).TextToDisplay = "a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.de-de.xlf"
$ws3.Hyperlinks.Item(8).Address = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/e6bf38c351c2af39466f1fd9fe06bae116651949/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/a8b88462-4116-4a82-a657-005344fb7aed.bd7b0ca6234c65b6273380aa27cd4258fac65a6e.de-de.xlf"
